$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new price text is numeric-looking (e.g. "380.34") must be
# forced to Text format first, otherwise Excel auto-converts them to numbers
# (losing the original inlineStr/text semantics and trailing zeros).
$textCells = @('D5','D6','D9','D10','D12','D14','D15','D17','D19','D20','D21','D23','D24','D26','D27','D30','D32','D33','D35','D39','D40','D43','D44','D45')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '51.615.49'
$ws.Range('E2').Value = '  +1.06%  '
$ws.Range('D3').Value = '2.984.47'
$ws.Range('E3').Value = '  +2.43%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '380.34'
$ws.Range('E5').Value = '  +2.96%  '
$ws.Range('D6').Value = '105.63'
$ws.Range('E6').Value = '  +1.76%  '
$ws.Range('E7').Value = '  +0.69%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = '0.595'
$ws.Range('E9').Value = '  +1.21%  '
$ws.Range('D10').Value = '37.44'
$ws.Range('E10').Value = '  +2.11%  '
$ws.Range('E11').Value = '  +0.35%  '
$ws.Range('D12').Value = '0.0846'
$ws.Range('E12').Value = '  +1.21%  '
$ws.Range('D13').Value = '3.453.79'
$ws.Range('E13').Value = '  +2.47%  '
$ws.Range('D14').Value = '18.47'
$ws.Range('E14').Value = '  +0.47%  '
$ws.Range('D15').Value = '7.57'
$ws.Range('E15').Value = '  +2.48%  '
$ws.Range('D16').Value = '2.984.13'
$ws.Range('E16').Value = '  +2.70%  '
$ws.Range('D17').Value = '0.972'
$ws.Range('E17').Value = '  +2.95%  '
$ws.Range('D18').Value = '51.547.74'
$ws.Range('E18').Value = '  +1.05%  '
$ws.Range('D19').Value = '3.35'
$ws.Range('E19').Value = '  +2.84%  '
$ws.Range('D20').Value = '7.44'
$ws.Range('E20').Value = '  +3.02%  '
$ws.Range('D21').Value = '13.00'
$ws.Range('E21').Value = '  +1.32%  '
$ws.Range('D22').Value = '0.0₃0965'
$ws.Range('E22').Value = '  +2.01%  '
$ws.Range('D23').Value = '69.55'
$ws.Range('E23').Value = '  +1.97%  '
$ws.Range('D24').Value = '262.34'
$ws.Range('E24').Value = '  +1.01%  '
$ws.Range('E25').Value = '  +5.93%  '
$ws.Range('D26').Value = '7.64'
$ws.Range('E26').Value = '  +24.97%  '
$ws.Range('D27').Value = '7.82'
$ws.Range('E27').Value = '  +12.31%  '
$ws.Range('E28').Value = '  +0.21%  '
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('D30').Value = '25.96'
$ws.Range('E30').Value = '  +0.97%  '
$ws.Range('E31').Value = '  +9.10%  '
$ws.Range('D32').Value = '9.95'
$ws.Range('E32').Value = '  +0.43%  '
$ws.Range('D33').Value = '35.47'
$ws.Range('E33').Value = '  +2.14%  '
$ws.Range('E34').Value = '  -1.93%  '
$ws.Range('D35').Value = '50.93'
$ws.Range('E35').Value = '  +0.19%  '
$ws.Range('E36').Value = '  +5.83%  '
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('E38').Value = '  +1.08%  '
$ws.Range('D39').Value = '17.28'
$ws.Range('E39').Value = '  +0.98%  '
$ws.Range('D40').Value = '2.63'
$ws.Range('E40').Value = '  -1.74%  '
$ws.Range('E41').Value = '  +0.66%  '
$ws.Range('E42').Value = '  +2.63%  '
$ws.Range('D43').Value = '125.20'
$ws.Range('E43').Value = '  +4.64%  '
$ws.Range('D44').Value = '21.91'
$ws.Range('E44').Value = '  -0.44%  '
$ws.Range('D45').Value = '0.293'
$ws.Range('E45').Value = '  +21.46%  '
$ws.Range('E46').Value = '  -1.62%  '
$ws.Range('E47').Value = '  +2.81%  '
$ws.Range('D48').Value = '2.049.31'
$ws.Range('E48').Value = '  +1.41%  '
$ws.Range('E49').Value = '  +1.81%  '
$ws.Range('E50').Value = '  +10.89%  '
$ws.Range('E51').Value = '  +2.68%  '
